$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.624.56"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.727.71"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.19"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.30"
$ws.Range("E6").Value = "  +6.55%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.550"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.726.14"
$ws.Range("E9").Value = "  +2.94%  "
$ws.Range("E10").Value = "  +2.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.367"
$ws.Range("E11").Value = "  +4.74%  "
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.78"
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.223.52"
$ws.Range("E15").Value = "  +3.18%  "
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.655.91"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.730.61"
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.95"
$ws.Range("E19").Value = "  +4.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "374.83"
$ws.Range("E20").Value = "  +4.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.69"
$ws.Range("E21").Value = "  +3.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.53"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.03"
$ws.Range("E23").Value = "  +5.73%  "
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.48"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.16"
$ws.Range("E27").Value = "  +3.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.866.41"
$ws.Range("E28").Value = "  +3.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000106"
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "594.27"
$ws.Range("E30").Value = "  +5.86%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +3.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.45"
$ws.Range("E33").Value = "  +4.17%  "
$ws.Range("E34").Value = "  +5.93%  "
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.92"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("E40").Value = "  +2.68%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.90"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.48"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.67"
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.97"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0312"
$ws.Range("E46").Value = "  -3.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.07"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("E48").Value = "  +5.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "155.43"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.94"
$ws.Range("E50").Value = "  +3.73%  "
$ws.Range("E51").Value = "  +5.89%  "
